$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species-observation data between row 3 and row 6 for the
# relevant columns (A, B, E, F, G, H, P, Q, R). All other columns are
# identical between the two rows so they are left untouched.
# Note: reads must use .Value2 (the .Value getter is unreliable here);
# writes use .Value.
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr6 = "$col`6"
    $val3 = $ws.Range($addr3).Value2
    $val6 = $ws.Range($addr6).Value2
    $ws.Range($addr3).Value = $val6
    $ws.Range($addr6).Value = $val3
}
